$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (2) .. L (12) are being reordered: the "Cash" column (currently L=12)
# moves to become the first data column, right after DATE (new B=2); the
# existing Salary..Party columns (old B..K = 2..11) all shift right by one
# (new C..L = 3..12). DATE (A/1), DAILY (M/13) and CUMULATIVE (N/14) do not move.
#
# new column index -> old column index
$newToOld = @{ 2 = 12; 3 = 2; 4 = 3; 5 = 4; 6 = 5; 7 = 6; 8 = 7; 9 = 8; 10 = 9; 11 = 10; 12 = 11 }

$lastRow = 63

for ($r = 1; $r -le $lastRow; $r++) {
    # Snapshot old B..L (2..12) contents for this row before overwriting anything.
    $old = @{}
    for ($c = 2; $c -le 12; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula()) {
            $old[$c] = @{ kind = "formula"; data = $cell.Formula() }
        } else {
            $old[$c] = @{ kind = "value"; data = $cell.Value() }
        }
    }

    # Write back in the new order.
    foreach ($newCol in $newToOld.Keys) {
        $oldCol = $newToOld[$newCol]
        $src = $old[$oldCol]
        $dst = $ws.Cells.Item($r, $newCol)
        if ($src.kind -eq "formula") {
            $dst.Formula = $src.data
        } else {
            $dst.Value = $src.data
        }
    }
}

# Column widths travel with the data: new B width = old L width, new C..L = old B..K.
# (Use .Width, which round-trips the stored width losslessly; .ColumnWidth
# rounds to 2 decimal places and would corrupt the value.)
$oldWidths = @{}
for ($c = 2; $c -le 12; $c++) {
    $oldWidths[$c] = $ws.Columns.Item($c).Width()
}
foreach ($newCol in $newToOld.Keys) {
    $oldCol = $newToOld[$newCol]
    $ws.Columns.Item($newCol).Width = $oldWidths[$oldCol]
}
